$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.746.33'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').Value = '3.407.72'
$ws.Range('E3').Value = '  -0.24%  '
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = '412.54'
$ws.Range('E5').Value = '  +0.80%  '
$ws.Range('D6').Value = '129.66'
$ws.Range('E6').Value = '  +0.35%  '
$ws.Range('D7').Value = '0.619'
$ws.Range('E7').Value = '  -3.05%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '0.723'
$ws.Range('E9').Value = '  -1.43%  '
$ws.Range('D10').Value = '0.136'
$ws.Range('E10').Value = '  -5.70%  '
$ws.Range('D11').Value = '42.59'
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('B12').Value = 'ShibaInu'
$ws.Range('C12').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D12').Value = '0.0000221'
$ws.Range('E12').Value = '  +0.66%  '
$ws.Range('D13').Value = '9.11'
$ws.Range('E13').Value = '  +2.11%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '3.943.04'
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('B15').Value = 'TRON'
$ws.Range('C15').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D15').Value = '0.140'
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').Value = '20.38'
$ws.Range('E16').Value = '  -2.19%  '
$ws.Range('D17').Value = '3.404.28'
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('D18').Value = '12.41'
$ws.Range('E18').Value = '  +2.51%  '
$ws.Range('D19').Value = '1.08'
$ws.Range('E19').Value = '  +0.99%  '
$ws.Range('D20').Value = '61.815.75'
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('D21').Value = '478.18'
$ws.Range('E21').Value = '  +16.47%  '
$ws.Range('D22').Value = '90.64'
$ws.Range('E22').Value = '  +1.41%  '
$ws.Range('D23').Value = '3.26'
$ws.Range('E23').Value = '  +2.66%  '
$ws.Range('D24').Value = '13.13'
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').Value = '3.29'
$ws.Range('E25').Value = '  +1.79%  '
$ws.Range('D26').Value = '9.71'
$ws.Range('E26').Value = '  +10.28%  '
$ws.Range('D27').Value = '33.08'
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').Value = '4.75'
$ws.Range('E28').Value = '  -0.87%  '
$ws.Range('E29').Value = '  +3.11%  '
$ws.Range('D30').Value = '11.88'
$ws.Range('E30').Value = '  -0.08%  '
$ws.Range('D31').Value = '2.64'
$ws.Range('E31').Value = '  -2.17%  '
$ws.Range('E32').Value = '  -1.58%  '
$ws.Range('D33').Value = '0.112'
$ws.Range('E33').Value = '  -3.65%  '
$ws.Range('D34').Value = '40.90'
$ws.Range('E34').Value = '  -4.70%  '
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  -0.71%  '
$ws.Range('E36').Value = '  +4.48%  '
$ws.Range('D37').Value = '0.0485'
$ws.Range('E37').Value = '  -2.94%  '
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').Value = '3.01'
$ws.Range('E39').Value = '  +3.24%  '
$ws.Range('D40').Value = '148.63'
$ws.Range('E40').Value = '  +5.16%  '
$ws.Range('D41').Value = '0.322'
$ws.Range('E41').Value = '  +3.55%  '
$ws.Range('D42').Value = '0.134'
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('D43').Value = '3.33'
$ws.Range('E43').Value = '  -0.76%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').Value = '2.06'
$ws.Range('E44').Value = '  +4.23%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').Value = '2.59'
$ws.Range('E45').Value = '  +7.29%  '
$ws.Range('D46').Value = '4.18'
$ws.Range('D47').Value = '2.36'
$ws.Range('E47').Value = '  +20.42%  '
$ws.Range('E48').Value = '  -1.14%  '
$ws.Range('D49').Value = '22.10'
$ws.Range('E49').Value = '  +0.75%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').Value = '112.68'
$ws.Range('E50').Value = '  +13.29%  '
$ws.Range('B51').Value = 'PEPE'
$ws.Range('C51').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D51').Value = '0.0₃0513'
$ws.Range('E51').Value = '  +14.57%  '

Write-Host "Applied 108 cell updates"
